$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48: convert A48 and D48 from text to true numeric values
$ws.Range("A48").Value = 27
$ws.Range("D48").Value = 20

# Row 49: new row of data (kept as text, matching the source import format)
$ws.Range("A49").Value = "26"
$ws.Range("B49").Value = "Partly Cloudy"
$ws.Range("C49").Value = "01/18/2025"
$ws.Range("D49").Value = "20"
